# "dev of alliance region"
# The "resources" sheet gains a couple of new alliance-related entries
# (STR_type / createAlliance / buyArchon), and the active tab/selection
# moves from the "rights" sheet back onto "resources".

$wb = $excel.ActiveWorkbook

$resources = $wb.Worksheets.Item("resources")

# Row 1 header: INT_level -> STR_type (INT_gem in B1 stays as-is)
$resources.Range("A1").Value = "STR_type"

# Row 2: the numeric "level" value becomes the "createAlliance" cost row
$resources.Range("A2").Value = "createAlliance"

# Row 3: new "buyArchon" cost row, with a price of 100
$resources.Range("A3").Value = "buyArchon"
$resources.Range("B3").Value = 100

# Bring "resources" to the front (tabSelected / activeTab) and move the
# selection to A4, matching where the author left off editing.
$resources.Activate()
$resources.Range("A4").Select() | Out-Null
